# Add a new "BREEDERSTATUS" choice-type block (pig breeder status) to the
# choice table on Sheet1, following the same id/parent_name/name/display_name
# layout already used by the other choice types (MEDICINE, VACCINE, BREED, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 28: header row for the new choice type itself (no parent_name)
$ws.Range("A28").Value = 27
$ws.Range("C28").Value = "BREEDERSTATUS"
$ws.Range("D28").Value = "สถานะการผสมพันธุ์"

# Rows 29-30: fill the name column first ...
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "BREEDERSTATUS"
$ws.Range("C29").Value = "BREEDERSTATUS_001"

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "BREEDERSTATUS"
$ws.Range("C30").Value = "BREEDERSTATUS_002"

# ... then the display_name column for those same two rows
$ws.Range("D29").Value = "ปกติ"
$ws.Range("D30").Value = "สำเร็จ"

# Row 31: entered as name + display_name together
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "BREEDERSTATUS"
$ws.Range("C31").Value = "BREEDERSTATUS_003"
$ws.Range("D31").Value = "ล้มเหลว"

# Match the formatting (wrap text, vertical-center) already used by the rest
# of the table, by copying it from the neighbouring rows instead of building
# a style from scratch.
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("C27:D27").Copy()
$ws.Range("C28:D28").PasteSpecial(-4122)
$ws.Range("A25:D27").Copy()
$ws.Range("A29:D31").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("E36").Select() | Out-Null
